# Commit: "added command line parsing"
#
# The underlying data change is that a new "AttackType" value of "none" is
# written into column F for every card row that previously had no
# AttackType set (i.e. summoners / support monsters with no melee, ranged
# or magic attack). This introduces one new shared string ("none") and
# populates column F for those 19 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (Card entries) whose AttackType (column F) needs to become "none".
$rows = @(2, 16, 17, 20, 23, 26, 33, 34, 53, 56, 63, 64, 72, 76, 79, 88, 91, 95, 96)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "none"
}

# Reflect the author's final cursor/selection position in the sheet view.
$ws.Range("G91").Select()
